$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (flow_base_flow): B4 114 -> 90 (kept as text, like the source data),
# D4 0.9 -> 0.71
$ws.Range("B4").Value = "'90"
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").Value = 0.71

# Row 7 (riparian): B7 0 -> 120 (kept as text), C7 stays "0" (text),
# D7 0 -> 0.94
$ws.Range("B7").Value = "'120"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'0"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = 0.94
